# Edit script for models/14sectors_cat/land.xlsx
# Commit message: "Updates on the cat model"
#
# Changes applied:
#  1. Update the comment on Europe!A15 (author "Iñigo:" -> "Enric:",
#     body "CAIT / http://www.wri.org/..." -> "EDGAR").
#  2. Update the historic CO2 emissions data row (Europe!C15:AA15 values
#     changed, Europe!AB15 filled in, and Europe!AC15:AJ15 appended).
#  3. Update the sheet view/selection state so that:
#       - Global!    selection -> I48:J48 (active cell J48)
#       - Europe!    selection -> AG15 (Europe becomes the active/tabSelected sheet)
#       - Catalonia! selection -> D3 (no longer the tabSelected sheet)
#     which also drives the workbook-level activeTab to point at Europe.

$wb = $excel.ActiveWorkbook
$wsGlobal    = $wb.Worksheets.Item("Global")
$wsEurope    = $wb.Worksheets.Item("Europe")
$wsCatalonia = $wb.Worksheets.Item("Catalonia")

# --- 1. Update the Europe!A15 comment text -------------------------------
$cmt = $wsEurope.Range("A15").Comment
$cmt.Text("Enric:`nEDGAR`n")

# --- 2. Update the historic CO2 emissions row (Europe row 15) ------------
$rowVals = New-Object 'object[,]' 1,34
$rowVals[0,0] = -0.34496931405827302
$rowVals[0,1] = -0.46814242262519801
$rowVals[0,2] = -0.46609129011263106
$rowVals[0,3] = -0.49515421907023205
$rowVals[0,4] = -0.43664204247300598
$rowVals[0,5] = -0.41879503143528102
$rowVals[0,6] = -0.42498349624917697
$rowVals[0,7] = -0.4181721856350456
$rowVals[0,8] = -0.38851950639383398
$rowVals[0,9] = -0.46841254940853805
$rowVals[0,10] = -0.33033722222970902
$rowVals[0,11] = -0.37792358760623901
$rowVals[0,12] = -0.38759988127086309
$rowVals[0,13] = -0.36132755199540101
$rowVals[0,14] = -0.37335515417665699
$rowVals[0,15] = -0.31262835196880601
$rowVals[0,16] = -0.36857917375828403
$rowVals[0,17] = -0.29009013615355
$rowVals[0,18] = -0.33518804786981399
$rowVals[0,19] = -0.399988360089993
$rowVals[0,20] = -0.31988691783600903
$rowVals[0,21] = -0.31755958479959595
$rowVals[0,22] = -0.31369644998465401
$rowVals[0,23] = -0.296045027887887
$rowVals[0,24] = -0.30862722177435603
$rowVals[0,25] = -0.309062138068257
$rowVals[0,26] = -0.29693685476449
$rowVals[0,27] = -0.27505925115152602
$rowVals[0,28] = -0.24473760610362399
$rowVals[0,29] = -0.24656475954628296
$rowVals[0,30] = -0.26476597050030798
$rowVals[0,31] = -0.22245570919959301
$rowVals[0,32] = -0.22419899523269299
$rowVals[0,33] = -0.21516157963982299
$wsEurope.Range("C15:AJ15").Value = $rowVals


# --- 3. Update selections / active sheet ----------------------------------
# Global: select I48:J48 with J48 as the (approximate) active cell.
$wsGlobal.Activate()
$wsGlobal.Range("I48:J48").Select()

# Catalonia: move the selection/active cell to D3 (matches the target,
# and Catalonia stops being the tab-selected sheet once Europe is
# activated below).
$wsCatalonia.Activate()
$wsCatalonia.Range("D3").Select()

# Europe: make it the active/tab-selected sheet with AG15 selected - this
# is the last sheet activated, so it also becomes the workbook's active tab.
$wsEurope.Activate()
$wsEurope.Range("AG15").Select()
